$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$data = @(
    @(44432, 3, 15, 85.37765382207297),
    @(44433, 2, 16, 91.06949741021117),
    @(44434, 2, 17, 96.76134099834937),
    @(44435, 4, 17, 96.76134099834937),
    @(44436, 2, 15, 85.37765382207297),
    @(44437, 0, 14, 79.68581023393477),
    @(44438, 5, 18, 102.4531845864876),
    @(44439, 0, 15, 85.37765382207297),
    @(44440, 0, 13, 73.99396664579658)
)

$startRow = 358
$endRow = $startRow + $data.Count - 1

# Copy formatting (style/number format/borders) from the last existing row down
# onto the new rows before writing values.
$srcRow = $ws.Range("A357:D357")
$destRows = $ws.Range("A" + $startRow + ":D" + $endRow)
$srcRow.Copy()
$destRows.PasteSpecial(-4122)

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "Added rows $startRow to $endRow"
